# Generate Report for Archive
#
# Update the localization-status report: the two files
# "0bc0b9d0-ab08-44ed-b618-4b7a032d88f9.md" and
# "51aff97c-a019-4e92-9f78-2813b4096071.md" have moved from the
# "Ready for handoff" status into the "In Translation" status.
# The third tracked file, "fc8c56c8-3c94-424e-a7be-ac152a5bd955.md",
# stays "Ready for handoff".

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns B (zh-cn) and C (de-de) hold the status ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("B4").Value = $newStatus
$overview.Range("C4").Value = $newStatus

# --- zh-cn sheet: column B holds the status ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("B4").Value = $newStatus

# --- de-de sheet: column B holds the status ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("B4").Value = $newStatus
